# Adding test Case to Search Module OPQA_1243
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# --- Update Runmode (col D) / Results (col E) for rows whose status changed ---

# Rows where the Results (col E) cell should now read FAIL instead of PASS
$ws.Cells.Item(23, 5).Value = "FAIL"
$ws.Cells.Item(24, 5).Value = "FAIL"
$ws.Cells.Item(25, 5).Value = "FAIL"
$ws.Cells.Item(30, 5).Value = "FAIL"

# New test cases added at the bottom of the sheet: Runmode flips from N to Y
$ws.Cells.Item(43, 4).Value = "Y"
$ws.Cells.Item(44, 4).Value = "Y"
$ws.Cells.Item(45, 4).Value = "Y"

# ... and their Results are filled in
$ws.Cells.Item(43, 5).Value = "FAIL"
$ws.Cells.Item(44, 5).Value = "PASS"
$ws.Cells.Item(45, 5).Value = "PASS"

# --- Column E width tweak ---
$ws.Columns.Item(5).ColumnWidth = 6.6

# --- Update the active view / selection to D44 ---
$ws.Activate() | Out-Null
$ws.Range("D44").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 41
$excel.ActiveWindow.ScrollColumn = 4
